$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows above the current row 150 ("Angeleno" 44258 record),
# shifting the existing rows 150-159 down to 152-161.
$ws.Rows.Item(150).Resize(2).Insert()

# Row 150: new weekly record - Ciruela / Larry Ann / Primera
$ws.Cells.Item(150, 1).Value = 3
$ws.Cells.Item(150, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(150, 3).Value = "Coquimbo"
$ws.Cells.Item(150, 4).Value = 44610
$ws.Cells.Item(150, 5).Value = 5
$ws.Cells.Item(150, 6).Value = "Fruta"
$ws.Cells.Item(150, 7).Value = 100103
$ws.Cells.Item(150, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(150, 9).Value = 100103002
$ws.Cells.Item(150, 10).Value = "Ciruela"
$ws.Cells.Item(150, 11).Value = "Larry Ann"
$ws.Cells.Item(150, 12).Value = "Primera"
$ws.Cells.Item(150, 13).Value = 45
$ws.Cells.Item(150, 14).Value = 14000
$ws.Cells.Item(150, 15).Value = 14000
$ws.Cells.Item(150, 16).Value = 14000
$ws.Cells.Item(150, 17).Value = "$/caja 15 kilos empedrada"
$ws.Cells.Item(150, 18).Value = "Provincia de San Felipe de Aconcagua"
$ws.Cells.Item(150, 19).Value = 933
$ws.Cells.Item(150, 20).Value = 15

# Row 151: new weekly record - Ciruela / Larry Ann / Segunda
$ws.Cells.Item(151, 1).Value = 3
$ws.Cells.Item(151, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(151, 3).Value = "Coquimbo"
$ws.Cells.Item(151, 4).Value = 44610
$ws.Cells.Item(151, 5).Value = 5
$ws.Cells.Item(151, 6).Value = "Fruta"
$ws.Cells.Item(151, 7).Value = 100103
$ws.Cells.Item(151, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(151, 9).Value = 100103002
$ws.Cells.Item(151, 10).Value = "Ciruela"
$ws.Cells.Item(151, 11).Value = "Larry Ann"
$ws.Cells.Item(151, 12).Value = "Segunda"
$ws.Cells.Item(151, 13).Value = 40
$ws.Cells.Item(151, 14).Value = 13000
$ws.Cells.Item(151, 15).Value = 13000
$ws.Cells.Item(151, 16).Value = 13000
$ws.Cells.Item(151, 17).Value = "$/caja 15 kilos empedrada"
$ws.Cells.Item(151, 18).Value = "Provincia de San Felipe de Aconcagua"
$ws.Cells.Item(151, 19).Value = 867
$ws.Cells.Item(151, 20).Value = 15
